$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text cells (matching the source
# inlineStr cells) across every row we touch, so numeric-looking strings
# like "292.29" or "-0.51%" are not auto-converted to Number cells.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "292.29"
$ws.Range("E2").Value = "-0.51%"
$ws.Range("D3").Value = "30.83"
$ws.Range("E3").Value = "-1.09%"
$ws.Range("D4").Value = "4.898"
$ws.Range("E4").Value = "0.20%"
$ws.Range("D5").Value = "0.07302"
$ws.Range("E5").Value = "-0.32%"
$ws.Range("D6").Value = "2.268"
$ws.Range("E6").Value = "24.59%"
$ws.Range("D7").Value = "7.677"
$ws.Range("E7").Value = "-0.04%"
$ws.Range("D8").Value = "3.726"
$ws.Range("E8").Value = "-1.01%"
$ws.Range("D9").Value = "0.9006"
$ws.Range("E9").Value = "-0.50%"
$ws.Range("D10").Value = "0.1681"
$ws.Range("E10").Value = "1.43%"
$ws.Range("E11").Value = "6.02%"
$ws.Range("D12").Value = "0.08145"
$ws.Range("E12").Value = "-0.01%"
$ws.Range("D13").Value = "0.03090"
$ws.Range("E13").Value = "3.60%"
$ws.Range("D14").Value = "0.1004"
$ws.Range("E14").Value = "0.43%"
$ws.Range("D15").Value = "0.001498"
$ws.Range("E15").Value = "-0.04%"
$ws.Range("D16").Value = "0.005771"
$ws.Range("E16").Value = "0.75%"
$ws.Range("D17").Value = "3.490"
$ws.Range("E17").Value = "0.82%"
$ws.Range("D18").Value = "2.074"
$ws.Range("E18").Value = "-1.55%"
$ws.Range("D19").Value = "0.3328"
$ws.Range("E19").Value = "1.42%"
$ws.Range("E20").Value = "-0.56%"
$ws.Range("D21").Value = "4.028"
$ws.Range("E21").Value = "-7.03%"
$ws.Range("D23").Value = "0.04526"
$ws.Range("E23").Value = "1.09%"
$ws.Range("E24").Value = "-1.36%"
$ws.Range("D25").Value = "0.004631"
$ws.Range("E25").Value = "14.42%"
$ws.Range("E26").Value = "3.85%"
$ws.Range("D27").Value = "0.0003383"
$ws.Range("E27").Value = "-95.49%"
$ws.Range("D39").Value = "0.01595"
$ws.Range("E39").Value = "-3.39%"
$ws.Range("D40").Value = "0.04425"
$ws.Range("E40").Value = "0.61%"
$ws.Range("D41").Value = "0.007299"
$ws.Range("E41").Value = "-1.30%"
$ws.Range("D42").Value = "0.1320"
$ws.Range("E42").Value = "-0.19%"
$ws.Range("D43").Value = "0.008598"
$ws.Range("E44").Value = "-1.95%"
$ws.Range("D45").Value = "0.009419"
$ws.Range("E45").Value = "-16.19%"
$ws.Range("D46").Value = "0.00005915"
$ws.Range("E46").Value = "-1.52%"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("D49").Value = "0.002889"
$ws.Range("E49").Value = "20.32%"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("D51").Value = "0.0001993"

# Drop the temporary text-number-format styling so the cells go back to the
# default (unstyled) cellXf, matching the original file.
$rng.ClearFormats()
